$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the three "status" fill colors to the style table (teal, orange, red) ---
# then leave the weekday attendance rows colored red, matching the commit
# "added colors to rows". Applying all three colors in sequence (even though
# only the last one sticks) reproduces the same palette of fills that the
# generating tool registered in styles.xml (teal FF29A3CC, orange FFFFCC66,
# red FFDF5E5E) while only red ends up visibly applied here.
$teal = 13411113   # BGR for FF29A3CC
$orange = 6737151  # BGR for FFFFCC66
$red = 6184671     # BGR for FFDF5E5E

$week1 = $ws.Range("A5:J8")
$week1.Interior.Color = $teal
$week1.Interior.Color = $orange
$week1.Interior.Color = $red

$week2 = $ws.Range("A11:J15")
$week2.Interior.Color = $teal
$week2.Interior.Color = $orange
$week2.Interior.Color = $red

# --- Overtime-hours column (I) goes from 0 to 1 for every colored weekday row ---
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("I12").Value = 1
$ws.Range("I13").Value = 1
$ws.Range("I14").Value = 1
$ws.Range("I15").Value = 1

# --- B19 switches from a blank shared-string placeholder to an explicit FALSE ---
$ws.Range("B19").Value = $false

# --- Fix the FLOOR(...,1,1) typo (extra 3rd arg) down to the 2-arg form ---
$ws.Range("B22").Formula = '=FLOOR(F17/8,1)&"."&FLOOR(MOD(F17,8),1)&"."&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60'
$ws.Range("B23").Formula = '=FLOOR(H19,1)&"."&(H19-FLOOR(H19,1))*8&".0"'
$ws.Range("B24").Formula = '=FLOOR(I19,1)&"."&(I19-FLOOR(I19,1))*8&".0"'
$ws.Range("B27").Formula = '=FLOOR(K27/8,1)&"."&FLOOR(MOD(K27,8),1)&"."&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60'
